$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Row 53 currently shows "Huesca" with Casos activos (C) = 0
# Row 54 currently shows "Huelva" with Casos activos (C) = 72
# The update swaps these two provinces (their shared-string order is swapped
# in the workbook) so that the row that used to read "Huesca" now reads
# "Huelva" with 72 active cases, and the row that used to read "Huelva" now
# reads "Huesca" with 0 active cases.
$ws.Range("A53").Value = "Huelva"
$ws.Range("C53").Value = 72

$ws.Range("A54").Value = "Huesca"
$ws.Range("C54").Value = 0

# Update the "last updated" timestamp string, shown in cell A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 10:16"
